$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after the existing "EC50" sheet so the
# sheet order / physical file numbering matches the target workbook.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "two_measurments"

# Pooled-SD-for-two-measurements sample data.
$ws2.Range("A1").Value = "a"
$ws2.Range("B1").Value = "b"
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = 3
$ws2.Range("A3").Value = 4
$ws2.Range("B3").Value = 5

$ws2.Range("B4").Select() | Out-Null
